# Remove two data rows from the "NEW" sheet:
#   - original row 7  (Caso 4662, ALTOLAGUIRRE 2397)
#   - original row 72 (Caso 4166, ALTOLAGUIRRE 2201)
# Deleting these rows shifts all subsequent rows up, which matches the
# target diff (dimension goes from A1:P78 to A1:P76).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NEW")

# Delete the higher-numbered row first so the lower row number (7) stays valid.
$ws.Rows.Item(72).Delete()
$ws.Rows.Item(7).Delete()
